# Add a "Eutrophication" worksheet (copy of "GWP") with new impact data, and
# add "price"/"price_unit" columns to the "info" sheet.

$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item("info")
$wsGWP  = $wb.Worksheets.Item("GWP")

# ---------------------------------------------------------------------------
# 1. Duplicate the GWP sheet (this gives an exact clone of values/styles,
#    including the theme-coloured fills that are otherwise hard to reproduce)
#    and place it after GWP (i.e. at the end of the workbook).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsGWP.Copy([System.Reflection.Missing]::Value, $lastSheet)
$wsEutro = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsEutro.Name = "Eutrophication"

# Remove the sort state that was copied along with the sheet.
try {
    $wsEutro.Sort.SortFields.Clear()
} catch {
}

# ---------------------------------------------------------------------------
# 2. Update the copied sheet's content for the Eutrophication indicator.
# ---------------------------------------------------------------------------
# Column B (unit) becomes "g N" instead of "kg CO2-eq" for every data row.
for ($r = 2; $r -le 12; $r++) {
    $wsEutro.Range("B" + $r).Value = "g N"
}

# Row 2 - Excavation
$wsEutro.Range("C2").Value = 7.3000000000000007
$wsEutro.Range("D2").Value = 7.1
$wsEutro.Range("E2").Value = 7.5

# Row 3 - Bricks
$wsEutro.Range("C3").Value = 4.8000000000000007
$wsEutro.Range("D3").Value = 4.5
$wsEutro.Range("E3").Value = 5.0999999999999996

# Row 4 - Cement
$wsEutro.Range("C4").Value = 12.8
$wsEutro.Range("D4").Value = 11.7
$wsEutro.Range("E4").Value = 13.899999999999999

# Row 5 - Concrete (no data)
$wsEutro.Range("C5").ClearContents()
$wsEutro.Range("D5").ClearContents()
$wsEutro.Range("E5").ClearContents()

# Row 6 - Gravel (no data)
$wsEutro.Range("C6").ClearContents()
$wsEutro.Range("D6").ClearContents()
$wsEutro.Range("E6").ClearContents()

# Row 7 - Plastic (no data)
$wsEutro.Range("C7").ClearContents()
$wsEutro.Range("D7").ClearContents()
$wsEutro.Range("E7").ClearContents()

# Row 8 - Sand (no data)
$wsEutro.Range("C8").ClearContents()
$wsEutro.Range("D8").ClearContents()
$wsEutro.Range("E8").ClearContents()

# Row 9 - StainlessSteel
$wsEutro.Range("C9").Value = 45.3
$wsEutro.Range("D9").Value = 32.700000000000003
$wsEutro.Range("E9").Value = 57

# Row 10 - StainlessSteelSheet
$wsEutro.Range("C10").Value = 8.5
$wsEutro.Range("D10").Value = 7.8
$wsEutro.Range("E10").Value = 9.1

# Row 11 - Steel
$wsEutro.Range("C11").Value = 27.5
$wsEutro.Range("D11").Value = 23.299999999999997
$wsEutro.Range("E11").Value = 33.5

# Row 12 - Wood
$wsEutro.Range("C12").Value = 1972
$wsEutro.Range("D12").Value = 1862
$wsEutro.Range("E12").Value = 2082

# ---------------------------------------------------------------------------
# 3. Add "price" / "price_unit" columns to the "info" sheet.
# ---------------------------------------------------------------------------
$wsInfo.Range("D1").Value = "price"
$wsInfo.Range("E1").Value = "price_unit"
$wsInfo.Range("A1").Copy()
$wsInfo.Range("D1:E1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Restore selections / active sheet to match the final state.
# ---------------------------------------------------------------------------
$wsEutro.Range("E26").Select()

$wsGWP.Range("D2").Select()

$wsInfo.Activate()
$wsInfo.Range("G12").Select()
